$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formulas first (before touching neighbouring cell number formats,
#     otherwise the newly-created shared formula picks up a neighbour's
#     number format) ---

# H column: shared formula across the whole block H27:H40
$ws.Range("H27:H40").Formula = "=F27+G27"

# J / K columns: extend the existing shared formulas down to the new rows
# (leave J33:J38 / K33:K38 untouched, add the new tail as its own block)
$ws.Range("J39:J40").Formula = "=G39-I39"
$ws.Range("K39:K40").NumberFormat = "0.0"
$ws.Range("K39:K40").Formula = "=J39/E39"

# --- New simulation results for 131072 cores (row 39) ---
$ws.Range("E39").Value = 118
$ws.Range("F39").Value = 74.9
$ws.Range("G39").Value = 340.3
$ws.Range("I39").NumberFormat = "0.00E+00"
$ws.Range("I39").Value = 0.02

# --- New simulation results for 262144 cores (row 40) ---
$ws.Range("E40").Value = 166
$ws.Range("F40").Value = 103.5
$ws.Range("G40").Value = 468.5
$ws.Range("I40").NumberFormat = "0.00E+00"
$ws.Range("I40").Value = 0.049

# --- sheet view: active selection ---
$ws.Range("K40").Select()
